# animal_sample_table_labeled_elements.xlsx edit script
# Implements:
#  - Removes the extra "test_animal_6" animal row (Animals!row7) and the
#    corresponding "test_animal_6_sample_1" sample row (Samples!row7), which
#    were invalid test rows.
#  - Replaces the glucose-[...] Infusate values (which included isotope
#    labels not actually present in glucose) with methionine-[...] values
#    on the Animals sheet (column E), since methionine supports all of the
#    currently-tested elements.
#  - Widens the Animals!E column a bit to fit the new text, and re-selects
#    row 7 (now blank) on the Animals and Samples sheets, matching the
#    selection state left behind in the authored workbook.

$wb = $excel.ActiveWorkbook

$wsAnimals = $wb.Worksheets.Item("Animals")
$wsSamples = $wb.Worksheets.Item("Samples")
$wsTreatments = $wb.Worksheets.Item("Treatments")
$wsTissues = $wb.Worksheets.Item("Tissues")

# --- Animals sheet: swap the isotope-labeled glucose values for methionine ---
$wsAnimals.Range("E2").Value = "methionine-[13C5]"
$wsAnimals.Range("E3").Value = "methionine-[15N1]"
$wsAnimals.Range("E4").Value = "methionine-[2H11]"
$wsAnimals.Range("E5").Value = "methionine-[17O2]"
$wsAnimals.Range("E6").Value = "methionine-[33S1]"

# --- Remove the bogus 6th test animal / sample rows ---
$wsAnimals.Rows.Item(7).Delete()
$wsSamples.Rows.Item(7).Delete()

# --- Widen the Infusate column on Animals to fit the new values ---
$wsAnimals.Range("E1").EntireColumn.ColumnWidth = 15.666666666666666

# --- Set explicit column widths that were (re)computed for the other sheets ---
$wsSamples.Range("A1").EntireColumn.ColumnWidth = 19.166666666666668
$wsSamples.Range("B1").EntireColumn.ColumnWidth = 11.833333333333334
$wsSamples.Range("C1").EntireColumn.ColumnWidth = 13.833333333333334
$wsSamples.Range("D1").EntireColumn.ColumnWidth = 14.5
$wsSamples.Range("E1").EntireColumn.ColumnWidth = 12.5
$wsSamples.Range("F1").EntireColumn.ColumnWidth = 11.0

$wsTreatments.Range("A1").EntireColumn.ColumnWidth = 14.333333333333334
$wsTreatments.Range("B1").EntireColumn.ColumnWidth = 129.66666666666666

$wsTissues.Range("A1").EntireColumn.ColumnWidth = 33.833333333333336
$wsTissues.Range("B1").EntireColumn.ColumnWidth = 59.333333333333336

# --- Restore the selections left in the sheets after editing ---
$wsAnimals.Rows.Item(7).Select()
$wsSamples.Rows.Item(7).Select()
$wsTreatments.Range("A1:B1").EntireColumn.Select()
$wsTissues.Range("B3").Select()

# Leave Animals as the active/visible sheet, as in the original workbook.
$wsAnimals.Activate()
